# SpoilsCardActions.xlsx -- "Finished excel sheet for expansion spoils"
# Rebuilds Sheet1 with a new "Bonus text" column (B) and "Restrictions" column (G),
# re-flowing the existing data and appending many new spoils cards (rows 12-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean so no stale cells from the old 6-column layout linger.
$ws.Cells.Clear()

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Bonus text"
$ws.Range("C1").Value = "Active Gains"
$ws.Range("D1").Value = "When Use Active"
$ws.Range("E1").Value = "#Uses"
$ws.Range("F1").Value = "Passive Gains"
$ws.Range("G1").Value = "Restrictions"

$ws.Range("A2").Value = "Jr. Fireside Camper's Guide"
$ws.Range("C2").Value = " +1 survival skill check success"
$ws.Range("D2").Value = "After failing survival skill check"
$ws.Range("E2").Value = "Once"
$ws.Range("F2").Value = " +1 psych resistence"

$ws.Range("A3").Value = "Lucky Playing Cards"
$ws.Range("C3").Value = "Play poker"
$ws.Range("D3").Value = "During Effects Phase"
$ws.Range("E3").Value = "Once Per Turn"
$ws.Range("F3").Value = "None"

$ws.Range("A4").Value = "Compass and Maps"
$ws.Range("C4").Value = "Ignore delays and Action cards causing unwanted movement of any kind"
$ws.Range("D4").Value = "Anytime"
$ws.Range("E4").Value = "Unlimited"
$ws.Range("F4").Value = " +1 movement"

$ws.Range("A5").Value = "Armored Humvee"
$ws.Range("B5").Value = "With 7.62mm M60A3 Turret"
$ws.Range("C5").Value = "None"
$ws.Range("D5").Value = "None"
$ws.Range("E5").Value = "None"
$ws.Range("F5").Value = " +2 movement; All hexes cost 1 movement"

$ws.Range("A6").Value = "John Rammbo's"
$ws.Range("B6").Value = "Compound Hunting Bow"
$ws.Range("C6").Value = "Choose 2 opposing characters and assign them each 1 point of damange"
$ws.Range("D6").Value = "End of PvP Round"
$ws.Range("E6").Value = "Once per PvP Round"
$ws.Range("F6").Value = "None"

$ws.Range("A7").Value = "Lucky Day"
$ws.Range("C7").Value = " +6 spoils cards"
$ws.Range("D7").Value = "Immediately"
$ws.Range("E7").Value = "Once"
$ws.Range("F7").Value = "None"
$ws.Range("G7").Value = "Can't start game with (All event cards I think)"

$ws.Range("A8").Value = "Cache of Combat Fatigues"
$ws.Range("C8").Value = "None"
$ws.Range("D8").Value = "None"
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "Can combine with Armor"

$ws.Range("A9").Value = "Some Jackasse's Goggles"
$ws.Range("C9").Value = "None"
$ws.Range("D9").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "None"

$ws.Range("A10").Value = "Industrial Welding & Cutting Torch"
$ws.Range("C10").Value = " +2 Salvage"
$ws.Range("D10").Value = "After successful Encounter/Mission"
$ws.Range("E10").Value = "Unlimited"
$ws.Range("F10").Value = "None"

$ws.Range("A11").Value = "5.56mm Light Machinegun"
$ws.Range("C11").Value = "None"
$ws.Range("D11").Value = "None"
$ws.Range("E11").Value = "None"
$ws.Range("F11").Value = " +1 prestige"

$ws.Range("A12").Value = "9mm Semi Automatic Pistol"
$ws.Range("C12").Value = "None"
$ws.Range("D12").Value = "None"
$ws.Range("E12").Value = "None"
$ws.Range("F12").Value = "None"

$ws.Range("A13").Value = "Five 9MM Semi Auto Pistols"
$ws.Range("C13").Value = "None"
$ws.Range("D13").Value = "None"
$ws.Range("E13").Value = "None"
$ws.Range("F13").Value = "None"

$ws.Range("A14").Value = "Armored Car"
$ws.Range("C14").Value = "Discard Ambush encounters for another from the same deck"
$ws.Range("D14").Value = "After drawing Ambush encounter"
$ws.Range("E14").Value = "Unlimited"
$ws.Range("F14").Value = " +2 movement"

$ws.Range("A15").Value = "7.62mm Machine Gun"
$ws.Range("C15").Value = "None"
$ws.Range("D15").Value = "None"
$ws.Range("E15").Value = "None"
$ws.Range("F15").Value = "None"

$ws.Range("A16").Value = "Kick Ass Sound System"
$ws.Range("C16").Value = "None"
$ws.Range("D16").Value = "None"
$ws.Range("E16").Value = "None"
$ws.Range("F16").Value = " +1 movement"
$ws.Range("G16").Value = "Can't have vehicle = bicycles or horses"

$ws.Range("A17").Value = "Jamison Bond's 77 Lotus"
$ws.Range("C17").Value = "4d6 damage to opponent OR 2d6 town health removal"
$ws.Range("D17").Value = "During PvP OR If within 1 hex of player's town"
$ws.Range("E17").Value = "Once OR Once"
$ws.Range("F17").Value = " +3 movment"

$ws.Range("A18").Value = "Pristine American Flag"
$ws.Range("C18").Value = " +2 salvage"
$ws.Range("D18").Value = "End Turn Phase"
$ws.Range("E18").Value = "Once per turn"
$ws.Range("F18").Value = " +2 prestige"

$ws.Range("A19").Value = "Forged Government Credentials"
$ws.Range("C19").Value = " +1 prestige"
$ws.Range("D19").Value = "Anytime"
$ws.Range("E19").Value = "Once"
$ws.Range("F19").Value = "None"

$ws.Range("A20").Value = "Gargantuan BBQ Grill"
$ws.Range("B20").Value = "With Cooler"
$ws.Range("C20").Value = "Discard Perishable encounters for another from the same deck"
$ws.Range("D20").Value = "After drawing Perishable encounter"
$ws.Range("E20").Value = "Unlimited"
$ws.Range("F20").Value = " -1 Movement"
$ws.Range("G20").Value = "4 wheeled vehicle"

$ws.Range("A21").Value = "9mm Auto Pistol"
$ws.Range("C21").Value = "None"
$ws.Range("D21").Value = "None"
$ws.Range("E21").Value = "None"
$ws.Range("F21").Value = "None"

# Column widths (new B/G columns added, existing columns resized/shifted).
$ws.Columns.Item(2).ColumnWidth = 24.666666666666668
$ws.Columns.Item(3).ColumnWidth = 66.0
$ws.Columns.Item(4).ColumnWidth = 41.666666666666664
$ws.Columns.Item(5).ColumnWidth = 20.0
$ws.Columns.Item(6).ColumnWidth = 38.166666666666664
$ws.Columns.Item(7).ColumnWidth = 34.333333333333336

# Selection moved to D7 in the saved file.
$ws.Range("D7").Select() | Out-Null

